$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New error-code rows (row 2 = code 0 "no error", row 3 = code 2 "origin return incomplete")
$ws.Range("C2").Value = 0
$ws.Range("D2").Value = "에러 없음"

$ws.Range("C3").Value = 2
$ws.Range("D3").Value = "원점 복귀 미완료"
$ws.Range("E3").Value = "원점 복귀"

# Scroll the window so column C is the left-most visible column (topLeftCell = C1)
$excel.ActiveWindow.ScrollColumn = 3
$excel.ActiveWindow.ScrollRow = 1

$ws.Range("E3").Select()
